# Extend the credentials test data with extra "Peter" rows (peter124..peter131),
# turning the single-example table into a two-dimensional data array for
# data-driven ("excel to two dim array dynamic") testing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startPassword = 124
$endPassword = 131
$startRow = 5

for ($pwd = $startPassword; $pwd -le $endPassword; $pwd++) {
    $row = $startRow + ($pwd - $startPassword)
    $ws.Cells.Item($row, 1).Value = "Peter"
    $ws.Cells.Item($row, 2).Value = "peter$pwd"
    $ws.Cells.Item($row, 3).Value = "Invalid credentials"
}

# Match the author's selection state after adding the rows.
$null = $ws.Range("A5:C12").Select()
